$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as text; force text format so
# Excel does not auto-convert numeric-looking strings (e.g. "324.15") to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.110.79'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.791.00'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.15'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4298'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.69'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07530'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.0000'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.76'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.169'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.351'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.777.67'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.91'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06352'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.28'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.964'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.118.66'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.44'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.151'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.38'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.39'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.988.68'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.02'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.170'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.729'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09004'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.510'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02331'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.104'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6474'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2121'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06067'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.184'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9994'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.892'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.62'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5998'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.710'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.55'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.993'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.154'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06957'

# Volume(1h) column (E) values.
$ws.Range("E2").Value = '  +1.58%  '
$ws.Range("E3").Value = '  +1.85%  '
$ws.Range("E4").Value = '  -0.52%  '
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("E7").Value = '  -3.38%  '
$ws.Range("E9").Value = '  -2.75%  '
$ws.Range("E10").Value = '  -3.17%  '
$ws.Range("E11").Value = '  -1.16%  '
$ws.Range("E12").Value = '  -0.38%  '
$ws.Range("E13").Value = '  -0.04%  '
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("E16").Value = '  +0.99%  '
$ws.Range("E17").Value = '  +0.74%  '
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("E19").Value = '  +1.65%  '
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("E21").Value = '  -0.81%  '
$ws.Range("E22").Value = '  -3.66%  '
$ws.Range("E23").Value = '  +1.46%  '
$ws.Range("E24").Value = '  -1.90%  '
$ws.Range("E25").Value = '  -8.14%  '
$ws.Range("E26").Value = '  +4.39%  '
$ws.Range("E27").Value = '  -2.04%  '
$ws.Range("E28").Value = '  +1.50%  '
$ws.Range("E29").Value = '  -7.27%  '
$ws.Range("E30").Value = '  -1.77%  '
$ws.Range("E31").Value = '  -3.37%  '
$ws.Range("E32").Value = '  -0.88%  '
$ws.Range("E33").Value = '  -2.89%  '
$ws.Range("E34").Value = '  -4.97%  '
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("E36").Value = '  -0.46%  '
$ws.Range("E37").Value = '  +0.17%  '
$ws.Range("E38").Value = '  -0.50%  '
$ws.Range("E39").Value = '  -3.19%  '
$ws.Range("E40").Value = '  -0.96%  '
$ws.Range("E41").Value = '  -0.72%  '
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("E43").Value = '  -0.34%  '
$ws.Range("E44").Value = '  -1.56%  '
$ws.Range("E45").Value = '  -1.29%  '
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("E47").Value = '  -1.26%  '
$ws.Range("E48").Value = '  -1.10%  '
$ws.Range("E49").Value = '  -0.49%  '
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("E51").Value = '  +0.79%  '
